$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = 23
$ws.Range("G24").Value = 1042.82
$ws.Range("B25").Value = 7162.22
$ws.Range("F42").Value = 130
$ws.Range("G42").Value = 7670
$ws.Range("B54").Value = 100594.78
$ws.Range("F67").Value = 19
$ws.Range("G67").Value = 638.4
$ws.Range("B69").Value = 40687.91
$ws.Range("F96").Value = 37
$ws.Range("G96").Value = 2346.91
$ws.Range("F100").Value = 170
$ws.Range("G100").Value = 23859.5
$ws.Range("B116").Value = 160896.15
$ws.Range("F149").Value = 17
$ws.Range("G149").Value = 3680.5
$ws.Range("F151").Value = 65
$ws.Range("G151").Value = 9323.6
$ws.Range("F152").Value = 52
$ws.Range("G152").Value = 4440.28
$ws.Range("F153").Value = 131
$ws.Range("G153").Value = 8554.299999999999
$ws.Range("F161").Value = 179
$ws.Range("G161").Value = 11892.76
$ws.Range("B166").Value = 125204.45
$ws.Range("F203").Value = 61
$ws.Range("G203").Value = 3265.33
$ws.Range("B205").Value = 4389.86
$ws.Range("F207").Value = 21
$ws.Range("G207").Value = 660.03
$ws.Range("F212").Value = 68
$ws.Range("G212").Value = 3161.32
$ws.Range("F213").Value = 52
$ws.Range("G213").Value = 4517.76
$ws.Range("F215").Value = 194
$ws.Range("G215").Value = 9019.059999999999
$ws.Range("F216").Value = 36
$ws.Range("G216").Value = 2705.04
$ws.Range("F217").Value = 50
$ws.Range("G217").Value = 352
$ws.Range("B221").Value = 48665.5
$ws.Range("F223").Value = 166
$ws.Range("G223").Value = 18973.8
$ws.Range("F224").Value = 2334
$ws.Range("G224").Value = 43179
$ws.Range("B229").Value = 68183.62
$ws.Range("F272").Value = 36
$ws.Range("G272").Value = 3855.24
$ws.Range("F279").Value = 141
$ws.Range("G279").Value = 16106.43
$ws.Range("F281").Value = 15
$ws.Range("G281").Value = 2879.4
$ws.Range("F283").Value = 118
$ws.Range("G283").Value = 17025.04
$ws.Range("F284").Value = 84
$ws.Range("G284").Value = 12052.32
$ws.Range("B325").Value = 173011.1
$ws.Range("F343").Value = 150
$ws.Range("G343").Value = 3363
$ws.Range("B348").Value = 32169.98
$ws.Range("F351").Value = 1
$ws.Range("G351").Value = 24.33
$ws.Range("F353").Value = 7
$ws.Range("G353").Value = 45.99
$ws.Range("B358").Value = 770.2
$ws.Range("F363").Value = 392
$ws.Range("G363").Value = 55111.28
$ws.Range("B365").Value = 70000
$ws.Range("F367").Value = 20
$ws.Range("G367").Value = 409.8
$ws.Range("B372").Value = 757.3099999999999
$ws.Range("B387").Value = 47097
$ws.Range("D387").Value = 112.28
$ws.Range("E387").Value = 134.16
$ws.Range("F387").Value = 15
$ws.Range("G387").Value = 1684.2
$ws.Range("B388").Value = 58047
$ws.Range("D388").Value = 105.54
$ws.Range("E388").Value = 126.1
$ws.Range("F388").Value = 55
$ws.Range("G388").Value = 5804.7
$ws.Range("F409").Value = 54
$ws.Range("G409").Value = 2559.6
$ws.Range("F410").Value = 114
$ws.Range("G410").Value = 4265.88
$ws.Range("F415").Value = 58
$ws.Range("G415").Value = 1928.5
$ws.Range("B421").Value = 106779.47
$ws.Range("F430").Value = 350
$ws.Range("G430").Value = 4602.5
$ws.Range("F431").Value = 461
$ws.Range("G431").Value = 5905.41
$ws.Range("F436").Value = 320
$ws.Range("G436").Value = 6313.6
$ws.Range("F438").Value = 370
$ws.Range("G438").Value = 6001.4
$ws.Range("F440").Value = 734
$ws.Range("G440").Value = 4829.72
$ws.Range("F446").Value = 427
$ws.Range("G446").Value = 6289.71
$ws.Range("B447").Value = 76295.24000000001
$ws.Range("F453").Value = 164
$ws.Range("G453").Value = 8213.120000000001
$ws.Range("F455").Value = 309
$ws.Range("G455").Value = 15474.72
$ws.Range("F457").Value = 9
$ws.Range("G457").Value = 1953.72
$ws.Range("F461").Value = 131
$ws.Range("G461").Value = 6341.71
$ws.Range("F463").Value = 319
$ws.Range("G463").Value = 3075.16
$ws.Range("B469").Value = 132540.22
$ws.Range("F487").Value = 42
$ws.Range("G487").Value = 6182.4
$ws.Range("B491").Value = 35919.88
$ws.Range("F496").Value = 455
$ws.Range("G496").Value = 9031.75
$ws.Range("B501").Value = 27605.11
$ws.Range("F523").Value = 34
$ws.Range("G523").Value = 1110.44
$ws.Range("B538").Value = 59710.01
$ws.Range("F545").Value = 78
$ws.Range("G545").Value = 7683
$ws.Range("B552").Value = 59783.73
$ws.Range("F570").Value = 15
$ws.Range("G570").Value = 1538.1
$ws.Range("F575").Value = 148
$ws.Range("G575").Value = 8122.24
$ws.Range("F579").Value = 36
$ws.Range("G579").Value = 4799.88
$ws.Range("F580").Value = 67
$ws.Range("G580").Value = 4927.85
$ws.Range("F581").Value = 113
$ws.Range("G581").Value = 7174.37
$ws.Range("F582").Value = 115
$ws.Range("G582").Value = 8030.45
$ws.Range("F586").Value = 154
$ws.Range("G586").Value = 13604.36
$ws.Range("B588").Value = 133675.11
$ws.Range("F599").Value = 1
$ws.Range("G599").Value = 3615.9
$ws.Range("B604").Value = 10869.14
$ws.Range("F632").Value = 11
$ws.Range("G632").Value = 6648.62
$ws.Range("F634").Value = 1
$ws.Range("G634").Value = 5524.18
$ws.Range("B635").Value = 12172.8
$ws.Range("F637").Value = 96
$ws.Range("G637").Value = 7957.44
$ws.Range("F638").Value = 114
$ws.Range("G638").Value = 13890.9
$ws.Range("F640").Value = 99
$ws.Range("G640").Value = 8470.440000000001
$ws.Range("F641").Value = 91
$ws.Range("G641").Value = 7785.96
$ws.Range("B645").Value = 115340.58
$ws.Range("F669").Value = 23
$ws.Range("G669").Value = 1875.88
$ws.Range("F672").Value = 25
$ws.Range("G672").Value = 3262.5
$ws.Range("F675").Value = 287
$ws.Range("G675").Value = 10699.36
$ws.Range("F679").Value = 38
$ws.Range("G679").Value = 4211.54
$ws.Range("F680").Value = 28
$ws.Range("G680").Value = 3780.28
$ws.Range("F681").Value = 151
$ws.Range("G681").Value = 18227.21
$ws.Range("B683").Value = 65640.52
$ws.Range("F688").Value = 41
$ws.Range("G688").Value = 4461.21
$ws.Range("F692").Value = 148
$ws.Range("G692").Value = 22260.68
$ws.Range("F695").Value = 38
$ws.Range("G695").Value = 2005.26
$ws.Range("F700").Value = 114
$ws.Range("G700").Value = 3433.68
$ws.Range("F704").Value = 69
$ws.Range("G704").Value = 2820.03
$ws.Range("B713").Value = 66952.63
$ws.Range("F751").Value = 2128
$ws.Range("G751").Value = 347098.08
$ws.Range("F752").Value = 222
$ws.Range("G752").Value = 62797.14
$ws.Range("F753").Value = 339
$ws.Range("G753").Value = 49036.35
$ws.Range("F754").Value = 50
$ws.Range("G754").Value = 1907
$ws.Range("F758").Value = 13
$ws.Range("G758").Value = 1671.54
$ws.Range("B759").Value = 492518.04
$ws.Range("B764").Value = 3034888.91
$ws.Range("B765").Value = 3034888.91
